# Apply "Add data for 2022-08-18" update to the carjacking-by-neighborhood
# pivot workbook: rename the "through Aug 09" sheet/label to "through Aug 10"
# and bump/add a handful of data points for 2022-08-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet tab and update its title text (B1 header + shared string)
$ws.Name = "Through 2022-08-10"
$ws.Range("B1").Value = "August 2022 (through August 10)"

# Cell value updates: column letter + row -> new value
$updates = @(
    @("J2",  4),
    @("R2",  3),
    @("J3",  1),
    @("AX3", 1),
    @("J4",  4),
    @("J5",  8),
    @("R5",  6),
    @("B7",  4),
    @("R7",  3),
    @("J8",  2),
    @("J12", 3),
    @("AP13",1),
    @("J17", 1),
    @("R21", 1),
    @("R27", 2),
    @("R29", 2),
    @("B35", 2),
    @("R38", 1),
    @("R39", 2),
    @("AX50",2),
    @("AX66",1),
    @("AX69",2),
    @("J74", 2),
    @("R95", 3),
    @("R97", 1)
)

foreach ($pair in $updates) {
    $ws.Range($pair[0]).Value = $pair[1]
}
